# Candidate power plants were assigned to the wrong producer / wrong rows.
# This script re-points the "times" and "scenario_data_emlab" scenario
# settings to the new run, and inserts the correct candidate power-plant
# rows (ahead of the existing entries) on the conventionals / renewables /
# storages / biogas sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "times": shift the simulation window forward by 4 years
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("times")
$ws.Range("B2").Value = 45291.99861111111
$ws.Range("B3").Value = 45656.99861111111

# ---------------------------------------------------------------
# Sheet "scenario_data_emlab": move scenario year 2020 -> 2024 and
# update commodity / CO2 prices accordingly
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("scenario_data_emlab")
$ws.Range("B1").Value = 2024
$ws.Range("B2").Value = 36.31999999999999
$ws.Range("B5").Value = 13.616
$ws.Range("B6").Value = 21.392
$ws.Range("B7").Value = 53.136

# ---------------------------------------------------------------
# Sheet "conventionals": insert the two candidate power plants that
# were missing, ahead of the existing (already built) plants, and fix
# up the efficiency values of the existing plants.
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("conventionals")
$ws.Rows("2:3").Insert()
$ws.Range("B2:G3").ClearFormats()
$ws.Range("A4").Copy()
$ws.Range("A2:A3").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 99991700006
$ws.Range("C2").Value = "NATURAL_GAS"
$ws.Range("D2").Value = 4.5
$ws.Range("E2").Value = 0.43
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 99990300008
$ws.Range("C3").Value = "NATURAL_GAS"
$ws.Range("D3").Value = 4.2
$ws.Range("E3").Value = 0.61
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1

# fix indices + efficiency of the existing (shifted) rows
$ws.Range("A4").Value = 2
$ws.Range("E4").Value = 0.61

$ws.Range("A5").Value = 3
$ws.Range("E5").Value = 0.33

$ws.Range("A6").Value = 4
$ws.Range("E6").Value = 0.35

$ws.Range("A7").Value = 5
$ws.Range("E7").Value = 0.33

$ws.Range("A8").Value = 6
$ws.Range("E8").Value = 0.33

$ws.Range("A9").Value = 7
$ws.Range("E9").Value = 0.43

# ---------------------------------------------------------------
# Sheet "renewables": insert the three candidate renewable plants
# that were missing, ahead of the existing (already built) plants.
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("renewables")
$ws.Rows("2:4").Insert()
$ws.Range("B2:I4").ClearFormats()
$ws.Range("A5").Copy()
$ws.Range("A2:A4").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 99992100002
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = "OtherPV"
$ws.Range("F2").Value = "-"
$ws.Range("G2").Value = "-"
$ws.Range("H2").Value = "-"
$ws.Range("I2").Value = "-"

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 99992400003
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 1.35
$ws.Range("E3").Value = "WindOn"
$ws.Range("F3").Value = "-"
$ws.Range("G3").Value = "-"
$ws.Range("H3").Value = "-"
$ws.Range("I3").Value = "-"

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 99992300007
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 2.7
$ws.Range("E4").Value = "WindOff"
$ws.Range("F4").Value = "-"
$ws.Range("G4").Value = "-"
$ws.Range("H4").Value = "-"
$ws.Range("I4").Value = "-"

# fix indices of the existing (shifted) rows
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4
$ws.Range("A7").Value = 5
$ws.Range("A8").Value = 6

# ---------------------------------------------------------------
# Sheet "storages": add the first (candidate) storage unit
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("storages")
$ws.Range("A1").Copy()
$ws.Range("A2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 99992600009
$ws.Range("C2").Value = "STORAGE"
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 1

# ---------------------------------------------------------------
# Sheet "biogas": insert the candidate biogas plant ahead of the
# existing (already built) one.
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("biogas")
$ws.Rows("2:2").Insert()
$ws.Range("B2:I2").ClearFormats()
$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 99990100004
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 1.9
$ws.Range("E2").Value = "Biogas"
$ws.Range("F2").Value = "-"
$ws.Range("G2").Value = "-"
$ws.Range("H2").Value = "-"
$ws.Range("I2").Value = "-"

# fix index of the existing (shifted) row
$ws.Range("A3").Value = 1
